$wb = $excel.ActiveWorkbook

# The workbook already has a "template" sheet (FTNC_Demand5) whose header
# row (B1:F1) and label cell (A2) carry the bold/bordered/centered style
# used throughout every FTNC_Demand5x sheet. Reuse that exact formatting
# for the new sheets instead of re-declaring it, so no new style/font gets
# introduced into styles.xml.
$template = $wb.Worksheets.Item(1)

# Data for the five new demand sheets (FTNC_Demand56 .. FTNC_Demand510),
# matching the "In-vehicle / At-stop / Extra / Tardiness / Total" layout
# used by the existing FTNC_Demand5x sheets.
$sheetsData = @(
    @{ Name = "FTNC_Demand56";  B2 = 4.421875;           C2 = 221.8636545138889;  D2 = 0;                 E2 = 0;                 F2 = 226.2855295138889 },
    @{ Name = "FTNC_Demand57";  B2 = 12.4249913434903;   C2 = 186.3187554948454;  D2 = 0;                 E2 = 0;                 F2 = 198.7437468383357 },
    @{ Name = "FTNC_Demand58";  B2 = 15.19801765927978;  C2 = 189.8765048023247;  D2 = 0;                 E2 = 37.42165927533734; F2 = 242.4961817369417 },
    @{ Name = "FTNC_Demand59";  B2 = 8.646684556786706;  C2 = 178.8261308411057;  D2 = 0;                 E2 = 4.956073777799302; F2 = 192.4288891756918 },
    @{ Name = "FTNC_Demand510"; B2 = 2285.421019712964;  C2 = 12672.64042375143;  D2 = 575.6307206712264; E2 = 0;                 F2 = 15533.6921641356 }
)

foreach ($sd in $sheetsData) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $sd.Name

    # Header row
    $ws.Range("B1").Value = "In-vehicle"
    $ws.Range("C1").Value = "At-stop"
    $ws.Range("D1").Value = "Extra"
    $ws.Range("E1").Value = "Tardiness"
    $ws.Range("F1").Value = "Total"

    # Data row
    $ws.Range("A2").Value = "FTNC"
    $ws.Range("B2").Value = $sd.B2
    $ws.Range("C2").Value = $sd.C2
    $ws.Range("D2").Value = $sd.D2
    $ws.Range("E2").Value = $sd.E2
    $ws.Range("F2").Value = $sd.F2

    # Copy the exact formatting (bold, thin border, centered/top aligned)
    # from the template sheet so the new sheets use the same style.
    $null = $template.Range("B1:F1").Copy()
    $null = $ws.Range("B1:F1").PasteSpecial(-4122)
    $null = $template.Range("A2").Copy()
    $null = $ws.Range("A2").PasteSpecial(-4122)

    $null = $ws.Range("A1").Select()
}

$null = $wb.Worksheets.Item(1).Select()
